$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "AutoTestAdmin@@AutoTestUser"
$ws.Range("C3").Value = "New Transmittal from Automation"
$ws.Range("D3").Value = "UnTick"
$ws.Range("E3").Value = "Correspondence"
$ws.Range("F3").Value = "Issued for Review"
$ws.Range("L3").Value = "Message for New transmittal"
$ws.Range("M3").Value = "Comments for Issued for Review"
